# Insert a new column B ("Source") before the existing "Description" column,
# shifting "Description" to column C, and populate the new column with the
# source label for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at B; existing column B ("Description") and
# everything to its right shifts one column to the right (becomes C).
$ws.Columns.Item(2).Insert()

# Determine the last used row (should be 99: header + 98 data rows).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Header for the new column.
$ws.Range("B1").Value = "Source"

# Fill the new column with the source text for every data row (2..lastRow).
$sourceRange = $ws.Range("B2:B$lastRow")
$sourceRange.Value = "Inventario IEET - Mamíferos"
